# Apply updated crypto price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="59.209.14"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Value = '  +0.55%  '

$ws.Range("D3").Formula = '="2.521.06"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Value = '  +0.90%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Formula = '="533.80"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  -0.76%  '

$ws.Range("D6").Formula = '="140.00"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  -2.52%  '

$ws.Range("E7").Value = '  +0.22%  '

$ws.Range("D8").Formula = '="0.565"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = '  -0.85%  '

$ws.Range("D9").Formula = '="2.526.00"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = '  +0.18%  '

$ws.Range("E10").Value = '  +0.41%  '

$ws.Range("E11").Value = '  +0.67%  '

$ws.Range("D12").Formula = '="5.50"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = '  -1.66%  '

$ws.Range("E13").Value = '  +2.10%  '

$ws.Range("D14").Formula = '="2.965.75"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = '  +0.98%  '

$ws.Range("D15").Formula = '="23.09"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = '  -1.09%  '

$ws.Range("D16").Formula = '="59.149.80"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = '  +0.61%  '

$ws.Range("D17").Formula = '="0.0000140"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = '  +0.76%  '

$ws.Range("D18").Formula = '="2.525.33"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = '  +0.47%  '

$ws.Range("D19").Formula = '="11.01"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  -1.47%  '

$ws.Range("D20").Formula = '="4.26"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = '  -0.02%  '

$ws.Range("D21").Formula = '="322.61"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  -0.28%  '

$ws.Range("E23").Value = '  +1.45%  '

$ws.Range("D24").Formula = '="62.22"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  +0.70%  '

$ws.Range("D25").Formula = '="0.425"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  -2.85%  '

$ws.Range("E26").Value = '  +1.79%  '

$ws.Range("D27").Formula = '="0.998"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  +0.35%  '

$ws.Range("D28").Formula = '="7.81"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = '  +0.77%  '

$ws.Range("D29").Formula = '="6.84"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  +2.68%  '

$ws.Range("D30").Formula = '="0.0₃0767"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = '  -0.44%  '

$ws.Range("D31").Formula = '="1.79"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = '  -0.10%  '

$ws.Range("D32").Formula = '="161.56"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = '  +2.97%  '

$ws.Range("E33").Value = '  +0.22%  '

$ws.Range("E34").Value = '  -5.34%  '

$ws.Range("E35").Value = '  +1.65%  '

$ws.Range("D36").Formula = '="18.49"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  -0.98%  '

$ws.Range("D37").Formula = '="4.24"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = '  -2.02%  '

$ws.Range("D38").Formula = '="1.59"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = '  -0.70%  '

$ws.Range("D39").Formula = '="36.95"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  +0.69%  '

$ws.Range("D40").Formula = '="3.65"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  -0.92%  '

$ws.Range("D41").Formula = '="0.806"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  -2.09%  '

$ws.Range("D42").Formula = '="5.24"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = '  -6.91%  '

$ws.Range("D43").Formula = '="280.47"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  -5.66%  '

$ws.Range("D44").Formula = '="0.998"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  +0.31%  '

$ws.Range("D45").Formula = '="10.88"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  +0.87%  '

$ws.Range("D46").Formula = '="0.596"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  -0.41%  '

$ws.Range("E47").Value = '  +0.51%  '

$ws.Range("D48").Formula = '="121.96"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  -1.37%  '

$ws.Range("D49").Formula = '="18.41"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = '  -0.41%  '

$ws.Range("D50").Formula = '="0.0512"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = '  -0.23%  '

$ws.Range("D51").Formula = '="0.0223"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = '  -1.85%  '

$excel.CutCopyMode = 0
